# Applies the dataset corrections described in the commit:
#  - updates several negativo/positivo (E/F) values across the existing rows
#  - removes the malformed duplicate row 99 ("Il Corrirere Della Sera" / Facebook / 1 / 0),
#    which shifts all subsequent rows (Strage di Cutro block) up by one and
#    shrinks the used range from A1:F111 down to A1:F110

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update negativo (E) / positivo (F) counts for existing rows ---
$ws.Range("E4").Value = 74
$ws.Range("F4").Value = 26

$ws.Range("E5").Value = 68
$ws.Range("F5").Value = 32

$ws.Range("E6").Value = 66
$ws.Range("F6").Value = 34

$ws.Range("E7").Value = 91
$ws.Range("F7").Value = 9

$ws.Range("E8").Value = 81

$ws.Range("E9").Value = 85

$ws.Range("E10").Value = 84
$ws.Range("F10").Value = 16

$ws.Range("E11").Value = 84

$ws.Range("E19").Value = 92

$ws.Range("E24").Value = 35
$ws.Range("F24").Value = 65

$ws.Range("E30").Value = 88
$ws.Range("F30").Value = 12

$ws.Range("E90").Value = 85

$ws.Range("E96").Value = 76

# --- Remove the bogus "Il Corrirere Della Sera" row (row 99) ---
# This shifts rows 100:111 up to 99:110, so the sheet ends up with
# a used range of A1:F110 instead of A1:F111.
$ws.Rows.Item(99).Delete()
